$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells D1, E1
$ws.Range("D1").Value = "hugh"
$ws.Range("E1").Value = "crime master gogo"

# Match the header formatting (bold, centered, bordered) used by C1
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

# Data rows for new columns D and E
$ws.Range("D2").Value = "Eligible"
$ws.Range("E2").Value = "Eligible"

$ws.Range("D3").Value = "Not Eligible"
$ws.Range("E3").Value = "Eligible"

$ws.Range("D4").Value = "Not Eligible"
$ws.Range("E4").Value = "Eligible"
